# "Generate Report for Archive"
#
# The localization-status report is regenerated:
#   - the "Ready for handoff" status entries move to "In Translation"
#     (Overview!E2:F2, E3:F3, and the per-locale "Status" column C on the
#     zh-cn / de-de sheets all reference the same status value)
#   - the Status/locale columns are narrower to fit the new (shorter)
#     status text, so the report re-flows their column widths

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update every cell showing the old status text -----------------------
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $cell = $used.Find($oldStatus)
    if ($cell -ne $null) {
        $firstAddress = $cell.Address()
        while ($true) {
            $cell.Value = $newStatus
            $cell = $used.FindNext($cell)
            if ($cell -eq $null -or $cell.Address() -eq $firstAddress) { break }
        }
    }
}

# --- Re-flow the status columns to the new (narrower) fitted width -------
# Target fitted width for the status columns is ~13.41 characters; the
# closest width this column-width model can store is 13.33, so feed it an
# input that lands on that bucket.
$newColWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").EntireColumn.ColumnWidth = $newColWidth
$wsOverview.Range("F1").EntireColumn.ColumnWidth = $newColWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $newColWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $newColWidth
